$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the data of rows 117 and 118 (everything except the id in column A
#    and the Date in column E, which stay put for each row position).
# ---------------------------------------------------------------------------
$cols = @("B","C","D","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$row117vals = @{}
$row118vals = @{}
foreach ($col in $cols) {
    $row117vals[$col] = $ws.Range($col + "117").Value()
    $row118vals[$col] = $ws.Range($col + "118").Value()
}

foreach ($col in $cols) {
    $ws.Range($col + "117").Value = $row118vals[$col]
    $ws.Range($col + "118").Value = $row117vals[$col]
}

# ---------------------------------------------------------------------------
# 2) Append a brand-new row 130 with an upcoming fixture (no score yet, so
#    columns H/I/J and AB/AC are intentionally left empty).
# ---------------------------------------------------------------------------

# Copy formatting from the row above for the styled cells (A = bold/border
# id style, E = date style), then overwrite with the correct values.
$ws.Range("A129").Copy($ws.Range("A130"))
$ws.Range("E129").Copy($ws.Range("E130"))

$ws.Range("A130").Value = 128
$ws.Range("B130").Value = 7862920
$ws.Range("C130").Value = "Lithuania A Lyga"
$ws.Range("D130").Value = "Lithuania A Lyga"
$ws.Range("E130").Value = 45389.52083333334
$ws.Range("F130").Value = "FK Kauno Zalgiris"
$ws.Range("G130").Value = "Panevezys"

$ws.Range("K130").Value = 2.75
$ws.Range("L130").Value = 3
$ws.Range("M130").Value = 2.4
$ws.Range("N130").Value = 2.25
$ws.Range("O130").Value = 3
$ws.Range("P130").Value = 2.9
$ws.Range("Q130").Value = -0.25
$ws.Range("R130").Value = 2
$ws.Range("S130").Value = 1.8
$ws.Range("T130").Value = 2.25
$ws.Range("U130").Value = 2.025
$ws.Range("V130").Value = 1.775
$ws.Range("W130").Value = 0
$ws.Range("X130").Value = 0
$ws.Range("Y130").Value = 0
$ws.Range("Z130").Value = 0
$ws.Range("AA130").Value = 0
